$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Global climate scheme (GCS))
$ws.Range("B2").Value = 55.3718303775116
$ws.Range("E2").Value = 56.7933697165468
$ws.Range("K2").Value = 54.8090076208069
$ws.Range("L2").Value = 49.0282605200155
$ws.Range("N2").Value = 47.4912333237318

# Row 3 (Belief about GCS support in own country)
$ws.Range("B3").Value = 40.5493307668479
$ws.Range("K3").Value = 36.9055302021312
$ws.Range("L3").Value = 36.9057907051389
$ws.Range("N3").Value = 41.6894250824717

# Row 4 (Belief about GCS support in the U.S.)
$ws.Range("B4").Value = 34.7872415482579
$ws.Range("K4").Value = 30.4357942185907
$ws.Range("N4").Value = 43.0029533260978
